# ----------------------------------------------------------------------------
# Error Calculations and Plots
#
# The first data row ("H 72") is removed from the missing-data sheet, shifting
# every subsequent record up by one row (dimension shrinks from A1:F63 to
# A1:F62). The regenerated missing-value mask also blanks/fills a handful of
# individual cells in the remaining rows (the underlying measurements are
# unchanged - only which cells are reported as missing differs).
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 2 ("H 72"); rows 3:63 shift up to become rows 2:62
$ws.Rows(2).Delete()

# Row 2 (H 1968 Sk 4): refresh missing-value mask
$ws.Range("D2").ClearContents()
$ws.Range("F2").ClearContents()

# Row 3 (H 12640 K XIII/3): refresh missing-value mask
$ws.Range("B3").Value = -19.8
$ws.Range("C3").Value = 11.7

# Row 4 (RM 2): refresh missing-value mask
$ws.Range("B4").ClearContents()

# Row 5 (RM 8): refresh missing-value mask
$ws.Range("C5").ClearContents()

# Row 6 (RM 9): refresh missing-value mask
$ws.Range("B6").Value = -18.7
$ws.Range("F6").Value = 0.71067

# Row 7 (RM 14): refresh missing-value mask
$ws.Range("B7").ClearContents()

# Row 10 (RM 32): refresh missing-value mask
$ws.Range("D10").Value = -13.8
$ws.Range("F10").ClearContents()

# Row 11 (RM 38): refresh missing-value mask
$ws.Range("F11").Value = 0.71122

# Row 13 (RM 52 a): refresh missing-value mask
$ws.Range("D13").ClearContents()

# Row 15 (RM 78): refresh missing-value mask
$ws.Range("B15").Value = -18.9
$ws.Range("F15").ClearContents()

# Row 16 (RM 81): refresh missing-value mask
$ws.Range("B16").ClearContents()

# Row 19 (RM 95): refresh missing-value mask
$ws.Range("D19").Value = -15.2

# Row 21 (RM 116): refresh missing-value mask
$ws.Range("D21").Value = -14.7

# Row 22 (RM 120): refresh missing-value mask
$ws.Range("D22").ClearContents()

# Row 24 (RM 134): refresh missing-value mask
$ws.Range("D24").ClearContents()

# Row 26 (RM 137): refresh missing-value mask
$ws.Range("D26").Value = -14.9

# Row 27 (RM 138): refresh missing-value mask
$ws.Range("B27").Value = -19.3
$ws.Range("D27").Value = -15.4

# Row 29 (RM 142a): refresh missing-value mask
$ws.Range("B29").ClearContents()
$ws.Range("D29").ClearContents()

# Row 30 (RM 142b): refresh missing-value mask
$ws.Range("D30").ClearContents()

# Row 31 (RM 145): refresh missing-value mask
$ws.Range("B31").Value = -19.5
$ws.Range("F31").Value = 0.7116

# Row 32 (RM 146): refresh missing-value mask
$ws.Range("B32").ClearContents()

# Row 34 (RM 158): refresh missing-value mask
$ws.Range("F34").Value = 0.70933

# Row 35 (RM 159): refresh missing-value mask
$ws.Range("F35").ClearContents()

# Row 37 (RM 167): refresh missing-value mask
$ws.Range("C37").Value = 12.1

# Row 38 (RM 170): refresh missing-value mask
$ws.Range("F38").ClearContents()

# Row 39 (RM 173): refresh missing-value mask
$ws.Range("B39").Value = -19.8
$ws.Range("C39").ClearContents()

# Row 40 (RM 178): refresh missing-value mask
$ws.Range("B40").ClearContents()

# Row 42 (RM 193): refresh missing-value mask
$ws.Range("D42").Value = -14.6

# Row 43 (RM 197): refresh missing-value mask
$ws.Range("D43").Value = -14.1

# Row 45 (RM 207): refresh missing-value mask
$ws.Range("C45").Value = 11.7
$ws.Range("D45").ClearContents()

# Row 46 (RM 208): refresh missing-value mask
$ws.Range("D46").ClearContents()

# Row 47 (RM 215): refresh missing-value mask
$ws.Range("C47").ClearContents()

# Row 51 (RM 233): refresh missing-value mask
$ws.Range("B51").Value = -20.5
$ws.Range("F51").Value = 0.7107599999999999

# Row 52 (SC 5): refresh missing-value mask
$ws.Range("B52").ClearContents()

# Row 53 (SC 66): refresh missing-value mask
$ws.Range("C53").Value = 10.5

# Row 55 (SC 101): refresh missing-value mask
$ws.Range("C55").ClearContents()
$ws.Range("D55").Value = -14.6
$ws.Range("F55").ClearContents()

# Row 56 (SC 103): refresh missing-value mask
$ws.Range("B56").Value = -19.2

# Row 57 (SC 105): refresh missing-value mask
$ws.Range("B57").ClearContents()

# Row 58 (SC 119): refresh missing-value mask
$ws.Range("D58").ClearContents()

# Row 59 (SC 120): refresh missing-value mask
$ws.Range("C59").Value = 11.4

# Row 60 (SC 132): refresh missing-value mask
$ws.Range("C60").Value = 15.3
$ws.Range("D60").Value = -13.7

# Row 61 (SC 193): refresh missing-value mask
$ws.Range("F61").Value = 0.71183

# Row 62 (SC 232): refresh missing-value mask
$ws.Range("C62").ClearContents()

Write-Output "Applied missing-data mask refresh; new used range $($ws.UsedRange.Address())"